$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column G (VMH_gene_IDs), shifting
# VMH_gene_IDs -> I1 and KEGG_genes -> J1
$ws.Range("G:H").Insert()

# Fill in the newly inserted headers
$ws.Range("G1").Value = "recon3_genes"
$ws.Range("H1").Value = "Lewis2010_genes"

# Match the style of the existing header row (bold, centered) for the new cells
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1:H1").HorizontalAlignment = -4108
